$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(83, 4).Value = 44895
$ws.Cells.Item(83, 10).Value = 2500
$ws.Cells.Item(83, 11).Value = 1500
$ws.Cells.Item(83, 12).Value = 2000
$ws.Cells.Item(83, 13).Value = 1750
$ws.Cells.Item(83, 16).Value = 1167
$ws.Cells.Item(84, 4).Value = 44497
$ws.Cells.Item(84, 10).Value = 3200
$ws.Cells.Item(84, 11).Value = 1300
$ws.Cells.Item(84, 12).Value = 1500
$ws.Cells.Item(84, 13).Value = 1400
$ws.Cells.Item(84, 16).Value = 933
$ws.Cells.Item(85, 4).Value = 44873
$ws.Cells.Item(85, 10).Value = 2400
$ws.Cells.Item(85, 11).Value = 1500
$ws.Cells.Item(85, 12).Value = 2000
$ws.Cells.Item(85, 13).Value = 1750
$ws.Cells.Item(85, 16).Value = 1167
$ws.Cells.Item(86, 4).Value = 44383
$ws.Cells.Item(86, 10).Value = 3200
$ws.Cells.Item(86, 11).Value = 1500
$ws.Cells.Item(86, 12).Value = 2000
$ws.Cells.Item(86, 13).Value = 1750
$ws.Cells.Item(86, 16).Value = 1167
$ws.Cells.Item(87, 4).Value = 44343
$ws.Cells.Item(87, 10).Value = 3340
$ws.Cells.Item(87, 11).Value = 1300
$ws.Cells.Item(87, 12).Value = 1500
$ws.Cells.Item(87, 13).Value = 1400
$ws.Cells.Item(87, 16).Value = 933
$ws.Cells.Item(88, 4).Value = 44474
$ws.Cells.Item(88, 10).Value = 2800
$ws.Cells.Item(88, 11).Value = 1500
$ws.Cells.Item(88, 12).Value = 2000
$ws.Cells.Item(88, 13).Value = 1750
$ws.Cells.Item(88, 16).Value = 1167
$ws.Cells.Item(89, 4).Value = 44418
$ws.Cells.Item(89, 10).Value = 3300
$ws.Cells.Item(89, 11).Value = 2000
$ws.Cells.Item(89, 12).Value = 2500
$ws.Cells.Item(89, 13).Value = 2250
$ws.Cells.Item(89, 16).Value = 1500
$ws.Cells.Item(90, 4).Value = 44231
$ws.Cells.Item(90, 10).Value = 2800
$ws.Cells.Item(90, 11).Value = 1300
$ws.Cells.Item(90, 12).Value = 1500
$ws.Cells.Item(90, 13).Value = 1400
$ws.Cells.Item(90, 16).Value = 933
$ws.Cells.Item(91, 4).Value = 44350
$ws.Cells.Item(91, 10).Value = 3340
$ws.Cells.Item(91, 11).Value = 1200
$ws.Cells.Item(91, 12).Value = 1500
$ws.Cells.Item(91, 13).Value = 1350
$ws.Cells.Item(91, 16).Value = 900
$ws.Cells.Item(92, 4).Value = 44341
$ws.Cells.Item(92, 10).Value = 3360
$ws.Cells.Item(92, 11).Value = 1300
$ws.Cells.Item(92, 12).Value = 1500
$ws.Cells.Item(92, 13).Value = 1400
$ws.Cells.Item(92, 16).Value = 933
$ws.Cells.Item(93, 4).Value = 44376
$ws.Cells.Item(93, 10).Value = 3200
$ws.Cells.Item(93, 11).Value = 1500
$ws.Cells.Item(93, 12).Value = 2000
$ws.Cells.Item(93, 13).Value = 1750
$ws.Cells.Item(93, 16).Value = 1167
$ws.Cells.Item(94, 4).Value = 44672
$ws.Cells.Item(94, 10).Value = 2400
$ws.Cells.Item(94, 11).Value = 2000
$ws.Cells.Item(94, 12).Value = 2500
$ws.Cells.Item(94, 13).Value = 2250
$ws.Cells.Item(94, 16).Value = 1500
$ws.Cells.Item(95, 4).Value = 44453
$ws.Cells.Item(95, 10).Value = 3200
$ws.Cells.Item(95, 11).Value = 2000
$ws.Cells.Item(95, 12).Value = 2500
$ws.Cells.Item(95, 13).Value = 2250
$ws.Cells.Item(95, 16).Value = 1500
$ws.Cells.Item(96, 4).Value = 44217
$ws.Cells.Item(96, 10).Value = 2800
$ws.Cells.Item(96, 11).Value = 1300
$ws.Cells.Item(96, 12).Value = 1500
$ws.Cells.Item(96, 13).Value = 1400
$ws.Cells.Item(96, 16).Value = 933
$ws.Cells.Item(97, 4).Value = 44663
$ws.Cells.Item(97, 10).Value = 2360
$ws.Cells.Item(97, 11).Value = 2000
$ws.Cells.Item(97, 12).Value = 2500
$ws.Cells.Item(97, 13).Value = 2250
$ws.Cells.Item(97, 16).Value = 1500
$ws.Cells.Item(98, 4).Value = 44540
$ws.Cells.Item(98, 10).Value = 3000
$ws.Cells.Item(98, 11).Value = 1500
$ws.Cells.Item(98, 12).Value = 2000
$ws.Cells.Item(98, 13).Value = 1750
$ws.Cells.Item(98, 16).Value = 1167
$ws.Cells.Item(99, 4).Value = 44420
$ws.Cells.Item(99, 10).Value = 3400
$ws.Cells.Item(99, 11).Value = 2000
$ws.Cells.Item(99, 12).Value = 2500
$ws.Cells.Item(99, 13).Value = 2250
$ws.Cells.Item(99, 16).Value = 1500
$ws.Cells.Item(100, 4).Value = 44224
$ws.Cells.Item(100, 10).Value = 2800
$ws.Cells.Item(100, 11).Value = 1300
$ws.Cells.Item(100, 12).Value = 1500
$ws.Cells.Item(100, 13).Value = 1400
$ws.Cells.Item(100, 16).Value = 933
$ws.Cells.Item(101, 4).Value = 44567
$ws.Cells.Item(101, 10).Value = 3200
$ws.Cells.Item(101, 11).Value = 3000
$ws.Cells.Item(101, 12).Value = 3500
$ws.Cells.Item(101, 13).Value = 3250
$ws.Cells.Item(101, 16).Value = 2167
$ws.Cells.Item(102, 4).Value = 44308
$ws.Cells.Item(102, 10).Value = 3200
$ws.Cells.Item(102, 11).Value = 1300
$ws.Cells.Item(102, 12).Value = 1500
$ws.Cells.Item(102, 13).Value = 1400
$ws.Cells.Item(102, 16).Value = 933
$ws.Cells.Item(103, 4).Value = 44348
$ws.Cells.Item(103, 10).Value = 3360
$ws.Cells.Item(103, 11).Value = 1300
$ws.Cells.Item(103, 12).Value = 1500
$ws.Cells.Item(103, 13).Value = 1400
$ws.Cells.Item(103, 16).Value = 933
$ws.Cells.Item(104, 4).Value = 44670
$ws.Cells.Item(104, 10).Value = 3000
$ws.Cells.Item(104, 11).Value = 2000
$ws.Cells.Item(104, 12).Value = 2500
$ws.Cells.Item(104, 13).Value = 2250
$ws.Cells.Item(104, 16).Value = 1500
$ws.Cells.Item(105, 4).Value = 44448
$ws.Cells.Item(105, 10).Value = 3200
$ws.Cells.Item(105, 11).Value = 2000
$ws.Cells.Item(105, 12).Value = 2500
$ws.Cells.Item(105, 13).Value = 2250
$ws.Cells.Item(105, 16).Value = 1500
$ws.Cells.Item(106, 4).Value = 44532
$ws.Cells.Item(106, 10).Value = 3260
$ws.Cells.Item(106, 11).Value = 1800
$ws.Cells.Item(106, 12).Value = 2000
$ws.Cells.Item(106, 13).Value = 1900
$ws.Cells.Item(106, 16).Value = 1267
$ws.Cells.Item(107, 4).Value = 44399
$ws.Cells.Item(107, 10).Value = 3320
$ws.Cells.Item(107, 11).Value = 1500
$ws.Cells.Item(107, 12).Value = 2000
$ws.Cells.Item(107, 13).Value = 1750
$ws.Cells.Item(107, 16).Value = 1167
$ws.Cells.Item(108, 4).Value = 44719
$ws.Cells.Item(108, 10).Value = 3200
$ws.Cells.Item(108, 11).Value = 1500
$ws.Cells.Item(108, 12).Value = 2000
$ws.Cells.Item(108, 13).Value = 1750
$ws.Cells.Item(108, 16).Value = 1167
$ws.Cells.Item(109, 4).Value = 44285
$ws.Cells.Item(109, 10).Value = 3400
$ws.Cells.Item(109, 11).Value = 2000
$ws.Cells.Item(109, 12).Value = 2500
$ws.Cells.Item(109, 13).Value = 2250
$ws.Cells.Item(109, 16).Value = 1500
$ws.Cells.Item(110, 4).Value = 44698
$ws.Cells.Item(110, 10).Value = 3200
$ws.Cells.Item(110, 11).Value = 1500
$ws.Cells.Item(110, 12).Value = 2000
$ws.Cells.Item(110, 13).Value = 1750
$ws.Cells.Item(110, 16).Value = 1167
$ws.Cells.Item(111, 4).Value = 44812
$ws.Cells.Item(111, 10).Value = 3000
$ws.Cells.Item(111, 11).Value = 2000
$ws.Cells.Item(111, 12).Value = 2500
$ws.Cells.Item(111, 13).Value = 2250
$ws.Cells.Item(111, 16).Value = 1500
$ws.Cells.Item(112, 4).Value = 44315
$ws.Cells.Item(112, 10).Value = 3120
$ws.Cells.Item(112, 11).Value = 1300
$ws.Cells.Item(112, 12).Value = 1500
$ws.Cells.Item(112, 13).Value = 1400
$ws.Cells.Item(112, 16).Value = 933
$ws.Cells.Item(113, 4).Value = 44357
$ws.Cells.Item(113, 10).Value = 3200
$ws.Cells.Item(113, 11).Value = 1300
$ws.Cells.Item(113, 12).Value = 1500
$ws.Cells.Item(113, 13).Value = 1400
$ws.Cells.Item(113, 16).Value = 933
$ws.Cells.Item(113, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(114, 4).Value = 44859
$ws.Cells.Item(114, 10).Value = 2400
$ws.Cells.Item(114, 11).Value = 1500
$ws.Cells.Item(114, 12).Value = 2000
$ws.Cells.Item(114, 13).Value = 1750
$ws.Cells.Item(114, 16).Value = 1167
$ws.Cells.Item(114, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(115, 4).Value = 44609
$ws.Cells.Item(115, 10).Value = 2600
$ws.Cells.Item(115, 11).Value = 2300
$ws.Cells.Item(115, 12).Value = 2500
$ws.Cells.Item(115, 13).Value = 2400
$ws.Cells.Item(115, 16).Value = 1600
$ws.Cells.Item(116, 4).Value = 44189
$ws.Cells.Item(116, 10).Value = 2600
$ws.Cells.Item(116, 11).Value = 1400
$ws.Cells.Item(116, 12).Value = 1500
$ws.Cells.Item(116, 13).Value = 1450
$ws.Cells.Item(116, 16).Value = 967
$ws.Cells.Item(117, 4).Value = 44819
$ws.Cells.Item(117, 10).Value = 2000
$ws.Cells.Item(117, 11).Value = 2000
$ws.Cells.Item(117, 12).Value = 2500
$ws.Cells.Item(117, 13).Value = 2250
$ws.Cells.Item(117, 16).Value = 1500
$ws.Cells.Item(118, 4).Value = 44229
$ws.Cells.Item(118, 10).Value = 3200
$ws.Cells.Item(118, 11).Value = 1300
$ws.Cells.Item(118, 12).Value = 1500
$ws.Cells.Item(118, 13).Value = 1400
$ws.Cells.Item(118, 16).Value = 933
$ws.Cells.Item(119, 4).Value = 44572
$ws.Cells.Item(119, 10).Value = 2900
$ws.Cells.Item(119, 11).Value = 3000
$ws.Cells.Item(119, 12).Value = 3500
$ws.Cells.Item(119, 13).Value = 3250
$ws.Cells.Item(119, 16).Value = 2167
$ws.Cells.Item(120, 4).Value = 44665
$ws.Cells.Item(120, 10).Value = 3000
$ws.Cells.Item(120, 11).Value = 2000
$ws.Cells.Item(120, 12).Value = 2500
$ws.Cells.Item(120, 13).Value = 2250
$ws.Cells.Item(120, 16).Value = 1500
$ws.Cells.Item(121, 4).Value = 44525
$ws.Cells.Item(121, 10).Value = 3200
$ws.Cells.Item(121, 11).Value = 1500
$ws.Cells.Item(121, 12).Value = 2000
$ws.Cells.Item(121, 13).Value = 1750
$ws.Cells.Item(121, 16).Value = 1167
$ws.Cells.Item(122, 4).Value = 44838
$ws.Cells.Item(122, 10).Value = 2800
$ws.Cells.Item(122, 11).Value = 2000
$ws.Cells.Item(122, 12).Value = 2500
$ws.Cells.Item(122, 13).Value = 2250
$ws.Cells.Item(122, 16).Value = 1500
$ws.Cells.Item(123, 4).Value = 44658
$ws.Cells.Item(123, 10).Value = 3000
$ws.Cells.Item(123, 11).Value = 2000
$ws.Cells.Item(123, 12).Value = 2500
$ws.Cells.Item(123, 13).Value = 2250
$ws.Cells.Item(123, 16).Value = 1500
$ws.Cells.Item(124, 4).Value = 44306
$ws.Cells.Item(124, 10).Value = 3400
$ws.Cells.Item(124, 11).Value = 2000
$ws.Cells.Item(124, 12).Value = 2500
$ws.Cells.Item(124, 13).Value = 2250
$ws.Cells.Item(124, 16).Value = 1500
$ws.Cells.Item(125, 4).Value = 44740
$ws.Cells.Item(125, 10).Value = 2460
$ws.Cells.Item(125, 11).Value = 1300
$ws.Cells.Item(125, 12).Value = 1500
$ws.Cells.Item(125, 13).Value = 1400
$ws.Cells.Item(125, 16).Value = 933
$ws.Cells.Item(126, 4).Value = 44462
$ws.Cells.Item(126, 10).Value = 3200
$ws.Cells.Item(126, 11).Value = 1500
$ws.Cells.Item(126, 12).Value = 2000
$ws.Cells.Item(126, 13).Value = 1750
$ws.Cells.Item(126, 16).Value = 1167
$ws.Cells.Item(127, 4).Value = 44313
$ws.Cells.Item(127, 10).Value = 3200
$ws.Cells.Item(127, 11).Value = 1300
$ws.Cells.Item(127, 12).Value = 1500
$ws.Cells.Item(127, 13).Value = 1400
$ws.Cells.Item(127, 16).Value = 933
$ws.Cells.Item(128, 4).Value = 44238
$ws.Cells.Item(128, 10).Value = 3200
$ws.Cells.Item(128, 11).Value = 1500
$ws.Cells.Item(128, 12).Value = 2000
$ws.Cells.Item(128, 13).Value = 1750
$ws.Cells.Item(128, 16).Value = 1167
$ws.Cells.Item(129, 4).Value = 44168
$ws.Cells.Item(129, 10).Value = 2800
$ws.Cells.Item(129, 11).Value = 1300
$ws.Cells.Item(129, 12).Value = 1500
$ws.Cells.Item(129, 13).Value = 1400
$ws.Cells.Item(129, 16).Value = 933
$ws.Cells.Item(130, 4).Value = 44467
$ws.Cells.Item(130, 10).Value = 3100
$ws.Cells.Item(130, 11).Value = 1500
$ws.Cells.Item(130, 12).Value = 2000
$ws.Cells.Item(130, 13).Value = 1750
$ws.Cells.Item(130, 16).Value = 1167
$ws.Cells.Item(131, 4).Value = 44537
$ws.Cells.Item(131, 10).Value = 3100
$ws.Cells.Item(131, 11).Value = 1500
$ws.Cells.Item(131, 12).Value = 2000
$ws.Cells.Item(131, 13).Value = 1750
$ws.Cells.Item(131, 16).Value = 1167
$ws.Cells.Item(132, 4).Value = 44784
$ws.Cells.Item(132, 10).Value = 2600
$ws.Cells.Item(132, 11).Value = 2000
$ws.Cells.Item(132, 12).Value = 2500
$ws.Cells.Item(132, 13).Value = 2250
$ws.Cells.Item(132, 16).Value = 1500
$ws.Cells.Item(133, 4).Value = 44245
$ws.Cells.Item(133, 10).Value = 3200
$ws.Cells.Item(133, 11).Value = 1500
$ws.Cells.Item(133, 12).Value = 2000
$ws.Cells.Item(133, 13).Value = 1750
$ws.Cells.Item(133, 16).Value = 1167
$ws.Cells.Item(134, 4).Value = 44411
$ws.Cells.Item(134, 10).Value = 3300
$ws.Cells.Item(134, 11).Value = 2000
$ws.Cells.Item(134, 12).Value = 2500
$ws.Cells.Item(134, 13).Value = 2250
$ws.Cells.Item(134, 16).Value = 1500
$ws.Cells.Item(135, 4).Value = 44791
$ws.Cells.Item(135, 10).Value = 2800
$ws.Cells.Item(135, 11).Value = 2000
$ws.Cells.Item(135, 12).Value = 2500
$ws.Cells.Item(135, 13).Value = 2250
$ws.Cells.Item(135, 16).Value = 1500
$ws.Cells.Item(136, 4).Value = 44397
$ws.Cells.Item(136, 10).Value = 3200
$ws.Cells.Item(136, 11).Value = 1500
$ws.Cells.Item(136, 12).Value = 2000
$ws.Cells.Item(136, 13).Value = 1750
$ws.Cells.Item(136, 16).Value = 1167
$ws.Cells.Item(137, 4).Value = 44775
$ws.Cells.Item(137, 10).Value = 2400
$ws.Cells.Item(137, 11).Value = 2000
$ws.Cells.Item(137, 12).Value = 2500
$ws.Cells.Item(137, 13).Value = 2250
$ws.Cells.Item(137, 16).Value = 1500
$ws.Cells.Item(138, 4).Value = 44544
$ws.Cells.Item(138, 10).Value = 2800
$ws.Cells.Item(138, 11).Value = 1500
$ws.Cells.Item(138, 12).Value = 2000
$ws.Cells.Item(138, 13).Value = 1750
$ws.Cells.Item(138, 16).Value = 1167
$ws.Cells.Item(139, 4).Value = 44868
$ws.Cells.Item(139, 10).Value = 2600
$ws.Cells.Item(139, 11).Value = 1500
$ws.Cells.Item(139, 12).Value = 2000
$ws.Cells.Item(139, 13).Value = 1750
$ws.Cells.Item(139, 16).Value = 1167
$ws.Cells.Item(140, 4).Value = 44336
$ws.Cells.Item(140, 10).Value = 3360
$ws.Cells.Item(140, 11).Value = 1300
$ws.Cells.Item(140, 12).Value = 1500
$ws.Cells.Item(140, 13).Value = 1400
$ws.Cells.Item(140, 16).Value = 933
$ws.Cells.Item(141, 4).Value = 44476
$ws.Cells.Item(141, 10).Value = 3000
$ws.Cells.Item(141, 11).Value = 1500
$ws.Cells.Item(141, 12).Value = 2000
$ws.Cells.Item(141, 13).Value = 1750
$ws.Cells.Item(141, 16).Value = 1167
$ws.Cells.Item(142, 4).Value = 44301
$ws.Cells.Item(142, 10).Value = 3200
$ws.Cells.Item(142, 11).Value = 2000
$ws.Cells.Item(142, 12).Value = 2500
$ws.Cells.Item(142, 13).Value = 2250
$ws.Cells.Item(142, 16).Value = 1500
$ws.Cells.Item(143, 4).Value = 44483
$ws.Cells.Item(143, 10).Value = 3300
$ws.Cells.Item(143, 11).Value = 1500
$ws.Cells.Item(143, 12).Value = 2000
$ws.Cells.Item(143, 13).Value = 1750
$ws.Cells.Item(143, 16).Value = 1167
$ws.Cells.Item(144, 4).Value = 44371
$ws.Cells.Item(144, 10).Value = 3300
$ws.Cells.Item(144, 11).Value = 1500
$ws.Cells.Item(144, 12).Value = 2000
$ws.Cells.Item(144, 13).Value = 1750
$ws.Cells.Item(144, 16).Value = 1167
$ws.Cells.Item(145, 4).Value = 44385
$ws.Cells.Item(145, 10).Value = 3320
$ws.Cells.Item(145, 11).Value = 1500
$ws.Cells.Item(145, 12).Value = 2000
$ws.Cells.Item(145, 13).Value = 1750
$ws.Cells.Item(145, 16).Value = 1167
$ws.Cells.Item(146, 4).Value = 44565
$ws.Cells.Item(146, 10).Value = 3000
$ws.Cells.Item(146, 11).Value = 3000
$ws.Cells.Item(146, 12).Value = 3500
$ws.Cells.Item(146, 13).Value = 3250
$ws.Cells.Item(146, 16).Value = 2167
$ws.Cells.Item(147, 4).Value = 44546
$ws.Cells.Item(147, 10).Value = 2800
$ws.Cells.Item(147, 11).Value = 2000
$ws.Cells.Item(147, 12).Value = 2500
$ws.Cells.Item(147, 13).Value = 2250
$ws.Cells.Item(147, 16).Value = 1500
$ws.Cells.Item(148, 4).Value = 44679
$ws.Cells.Item(148, 10).Value = 3000
$ws.Cells.Item(148, 11).Value = 2500
$ws.Cells.Item(148, 12).Value = 3000
$ws.Cells.Item(148, 13).Value = 2750
$ws.Cells.Item(148, 16).Value = 1833
$ws.Cells.Item(149, 4).Value = 44614
$ws.Cells.Item(149, 10).Value = 2300
$ws.Cells.Item(149, 11).Value = 2300
$ws.Cells.Item(149, 12).Value = 2500
$ws.Cells.Item(149, 13).Value = 2400
$ws.Cells.Item(149, 16).Value = 1600
$ws.Cells.Item(150, 4).Value = 44847
$ws.Cells.Item(150, 10).Value = 2800
$ws.Cells.Item(150, 11).Value = 2000
$ws.Cells.Item(150, 12).Value = 2500
$ws.Cells.Item(150, 13).Value = 2250
$ws.Cells.Item(150, 16).Value = 1500
$ws.Cells.Item(151, 4).Value = 44742
$ws.Cells.Item(151, 10).Value = 3000
$ws.Cells.Item(151, 11).Value = 1300
$ws.Cells.Item(151, 12).Value = 1500
$ws.Cells.Item(151, 13).Value = 1400
$ws.Cells.Item(151, 16).Value = 933
$ws.Cells.Item(152, 4).Value = 44187
$ws.Cells.Item(152, 10).Value = 3100
$ws.Cells.Item(152, 11).Value = 1400
$ws.Cells.Item(152, 12).Value = 1500
$ws.Cells.Item(152, 13).Value = 1450
$ws.Cells.Item(152, 16).Value = 967
$ws.Cells.Item(153, 4).Value = 44628
$ws.Cells.Item(153, 10).Value = 2400
$ws.Cells.Item(153, 11).Value = 2500
$ws.Cells.Item(153, 12).Value = 3000
$ws.Cells.Item(153, 13).Value = 2750
$ws.Cells.Item(153, 16).Value = 1833
$ws.Cells.Item(154, 4).Value = 44782
$ws.Cells.Item(154, 10).Value = 2600
$ws.Cells.Item(154, 11).Value = 2000
$ws.Cells.Item(154, 12).Value = 2500
$ws.Cells.Item(154, 13).Value = 2250
$ws.Cells.Item(154, 16).Value = 1500
$ws.Cells.Item(155, 4).Value = 44441
$ws.Cells.Item(155, 10).Value = 3200
$ws.Cells.Item(155, 11).Value = 2000
$ws.Cells.Item(155, 12).Value = 2500
$ws.Cells.Item(155, 13).Value = 2250
$ws.Cells.Item(155, 16).Value = 1500
$ws.Cells.Item(156, 4).Value = 44558
$ws.Cells.Item(156, 10).Value = 2900
$ws.Cells.Item(156, 11).Value = 2300
$ws.Cells.Item(156, 12).Value = 2500
$ws.Cells.Item(156, 13).Value = 2400
$ws.Cells.Item(156, 16).Value = 1600
$ws.Cells.Item(157, 4).Value = 44649
$ws.Cells.Item(157, 10).Value = 2800
$ws.Cells.Item(157, 11).Value = 2300
$ws.Cells.Item(157, 12).Value = 2500
$ws.Cells.Item(157, 13).Value = 2400
$ws.Cells.Item(157, 16).Value = 1600
$ws.Cells.Item(158, 4).Value = 44196
$ws.Cells.Item(158, 10).Value = 3200
$ws.Cells.Item(158, 11).Value = 1400
$ws.Cells.Item(158, 12).Value = 1500
$ws.Cells.Item(158, 13).Value = 1450
$ws.Cells.Item(158, 16).Value = 967
$ws.Cells.Item(159, 4).Value = 44425
$ws.Cells.Item(159, 10).Value = 3200
$ws.Cells.Item(159, 11).Value = 2000
$ws.Cells.Item(159, 12).Value = 2500
$ws.Cells.Item(159, 13).Value = 2250
$ws.Cells.Item(159, 16).Value = 1500
$ws.Cells.Item(160, 4).Value = 44581
$ws.Cells.Item(160, 10).Value = 3100
$ws.Cells.Item(160, 11).Value = 2500
$ws.Cells.Item(160, 12).Value = 3000
$ws.Cells.Item(160, 13).Value = 2750
$ws.Cells.Item(160, 16).Value = 1833
$ws.Cells.Item(161, 4).Value = 44406
$ws.Cells.Item(161, 10).Value = 3400
$ws.Cells.Item(161, 11).Value = 2000
$ws.Cells.Item(161, 12).Value = 2500
$ws.Cells.Item(161, 13).Value = 2250
$ws.Cells.Item(161, 16).Value = 1500
$ws.Cells.Item(162, 4).Value = 44161
$ws.Cells.Item(162, 10).Value = 3100
$ws.Cells.Item(162, 11).Value = 1300
$ws.Cells.Item(162, 12).Value = 1500
$ws.Cells.Item(162, 13).Value = 1400
$ws.Cells.Item(162, 16).Value = 933
$ws.Cells.Item(163, 4).Value = 44446
$ws.Cells.Item(163, 10).Value = 3200
$ws.Cells.Item(163, 11).Value = 2000
$ws.Cells.Item(163, 12).Value = 2500
$ws.Cells.Item(163, 13).Value = 2250
$ws.Cells.Item(163, 16).Value = 1500
$ws.Cells.Item(164, 4).Value = 44334
$ws.Cells.Item(164, 10).Value = 3440
$ws.Cells.Item(164, 11).Value = 1300
$ws.Cells.Item(164, 12).Value = 1500
$ws.Cells.Item(164, 13).Value = 1400
$ws.Cells.Item(164, 16).Value = 933
$ws.Cells.Item(165, 4).Value = 44488
$ws.Cells.Item(165, 10).Value = 3000
$ws.Cells.Item(165, 11).Value = 1300
$ws.Cells.Item(165, 12).Value = 1500
$ws.Cells.Item(165, 13).Value = 1400
$ws.Cells.Item(165, 16).Value = 933
$ws.Cells.Item(166, 4).Value = 44243
$ws.Cells.Item(166, 10).Value = 3200
$ws.Cells.Item(166, 11).Value = 1500
$ws.Cells.Item(166, 12).Value = 2000
$ws.Cells.Item(166, 13).Value = 1750
$ws.Cells.Item(166, 16).Value = 1167
$ws.Cells.Item(167, 4).Value = 44294
$ws.Cells.Item(167, 10).Value = 3000
$ws.Cells.Item(167, 11).Value = 2000
$ws.Cells.Item(167, 12).Value = 2500
$ws.Cells.Item(167, 13).Value = 2250
$ws.Cells.Item(167, 16).Value = 1500
$ws.Cells.Item(168, 4).Value = 44413
$ws.Cells.Item(168, 10).Value = 3360
$ws.Cells.Item(168, 11).Value = 2000
$ws.Cells.Item(168, 12).Value = 2500
$ws.Cells.Item(168, 13).Value = 2250
$ws.Cells.Item(168, 16).Value = 1500
$ws.Cells.Item(169, 4).Value = 44630
$ws.Cells.Item(169, 10).Value = 2000
$ws.Cells.Item(169, 11).Value = 2500
$ws.Cells.Item(169, 12).Value = 3000
$ws.Cells.Item(169, 13).Value = 2750
$ws.Cells.Item(169, 16).Value = 1833
$ws.Cells.Item(170, 4).Value = 44656
$ws.Cells.Item(170, 10).Value = 2400
$ws.Cells.Item(170, 11).Value = 2000
$ws.Cells.Item(170, 12).Value = 2500
$ws.Cells.Item(170, 13).Value = 2250
$ws.Cells.Item(170, 16).Value = 1500
$ws.Cells.Item(171, 4).Value = 44810
$ws.Cells.Item(171, 10).Value = 2400
$ws.Cells.Item(171, 11).Value = 2000
$ws.Cells.Item(171, 12).Value = 2500
$ws.Cells.Item(171, 13).Value = 2250
$ws.Cells.Item(171, 16).Value = 1500
$ws.Cells.Item(172, 4).Value = 44595
$ws.Cells.Item(172, 10).Value = 3000
$ws.Cells.Item(172, 11).Value = 2500
$ws.Cells.Item(172, 12).Value = 2800
$ws.Cells.Item(172, 13).Value = 2650
$ws.Cells.Item(172, 16).Value = 1767
$ws.Cells.Item(173, 4).Value = 44203
$ws.Cells.Item(173, 10).Value = 2800
$ws.Cells.Item(173, 11).Value = 1300
$ws.Cells.Item(173, 12).Value = 1500
$ws.Cells.Item(173, 13).Value = 1400
$ws.Cells.Item(173, 16).Value = 933
$ws.Cells.Item(174, 4).Value = 44278
$ws.Cells.Item(174, 10).Value = 3400
$ws.Cells.Item(174, 11).Value = 2000
$ws.Cells.Item(174, 12).Value = 2500
$ws.Cells.Item(174, 13).Value = 2250
$ws.Cells.Item(174, 16).Value = 1500
$ws.Cells.Item(175, 4).Value = 44495
$ws.Cells.Item(175, 10).Value = 2860
$ws.Cells.Item(175, 11).Value = 1300
$ws.Cells.Item(175, 12).Value = 1500
$ws.Cells.Item(175, 13).Value = 1400
$ws.Cells.Item(175, 16).Value = 933
$ws.Cells.Item(176, 4).Value = 44763
$ws.Cells.Item(176, 10).Value = 2000
$ws.Cells.Item(176, 11).Value = 2000
$ws.Cells.Item(176, 12).Value = 2500
$ws.Cells.Item(176, 13).Value = 2250
$ws.Cells.Item(176, 16).Value = 1500
$ws.Cells.Item(177, 4).Value = 44455
$ws.Cells.Item(177, 10).Value = 3200
$ws.Cells.Item(177, 11).Value = 2000
$ws.Cells.Item(177, 12).Value = 2500
$ws.Cells.Item(177, 13).Value = 2250
$ws.Cells.Item(177, 16).Value = 1500
$ws.Cells.Item(178, 4).Value = 44299
$ws.Cells.Item(178, 10).Value = 3400
$ws.Cells.Item(178, 11).Value = 2000
$ws.Cells.Item(178, 12).Value = 2500
$ws.Cells.Item(178, 13).Value = 2250
$ws.Cells.Item(178, 16).Value = 1500
$ws.Cells.Item(179, 4).Value = 44754
$ws.Cells.Item(179, 10).Value = 2400
$ws.Cells.Item(179, 11).Value = 1500
$ws.Cells.Item(179, 12).Value = 2000
$ws.Cells.Item(179, 13).Value = 1750
$ws.Cells.Item(179, 16).Value = 1167
$ws.Cells.Item(180, 4).Value = 44252
$ws.Cells.Item(180, 10).Value = 3600
$ws.Cells.Item(180, 11).Value = 1500
$ws.Cells.Item(180, 12).Value = 2000
$ws.Cells.Item(180, 13).Value = 1750
$ws.Cells.Item(180, 16).Value = 1167
$ws.Cells.Item(181, 4).Value = 44364
$ws.Cells.Item(181, 10).Value = 3200
$ws.Cells.Item(181, 11).Value = 1500
$ws.Cells.Item(181, 12).Value = 2000
$ws.Cells.Item(181, 13).Value = 1750
$ws.Cells.Item(181, 16).Value = 1167
$ws.Cells.Item(182, 4).Value = 44691
$ws.Cells.Item(182, 10).Value = 3400
$ws.Cells.Item(182, 11).Value = 1500
$ws.Cells.Item(182, 12).Value = 2000
$ws.Cells.Item(182, 13).Value = 1750
$ws.Cells.Item(182, 16).Value = 1167
$ws.Cells.Item(183, 4).Value = 44714
$ws.Cells.Item(183, 10).Value = 3200
$ws.Cells.Item(183, 11).Value = 1500
$ws.Cells.Item(183, 12).Value = 2000
$ws.Cells.Item(183, 13).Value = 1750
$ws.Cells.Item(183, 16).Value = 1167

# New row 184 (copy of former row 183)
$ws.Cells.Item(184, 1).Value = 8
$ws.Cells.Item(184, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 44236
$ws.Cells.Item(184, 4).NumberFormat = $ws.Cells.Item(183, 4).NumberFormat
$ws.Cells.Item(184, 5).Value = 4
$ws.Cells.Item(184, 6).Value = 100112044
$ws.Cells.Item(184, 7).Value = "Perejil"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 3200
$ws.Cells.Item(184, 11).Value = 1500
$ws.Cells.Item(184, 12).Value = 2000
$ws.Cells.Item(184, 13).Value = 1750
$ws.Cells.Item(184, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(184, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(184, 16).Value = 1167
$ws.Cells.Item(184, 17).Value = 1.5
$ws.Cells.Item(184, 18).Value = "Hortaliza"